# Remove obsolete query columns (mort, actualCountID, markCode) from the
# Recaptures-EDIQuery raw-data extract, shifting the remaining columns left.
# Deleting right-to-left keeps the earlier column letters valid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("W:W").Delete()
$ws.Columns("O:O").Delete()
$ws.Columns("J:J").Delete()
